# Add a new "2.3.1" benchmark row (httk v2.3.1 - "Minor bug fixes") to Sheet1,
# grow Table1 to cover it, and move the selection down to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 25

$ws.Cells.Item($row, 1).Value  = "2.3.1"
$ws.Cells.Item($row, 2).Value  = 1023
$ws.Cells.Item($row, 3).Value  = 0.99990000000000001
$ws.Cells.Item($row, 4).Value  = 1
$ws.Cells.Item($row, 5).Value  = 1
$ws.Cells.Item($row, 6).Value  = 1.0109999999999999
$ws.Cells.Item($row, 7).Value  = 352
$ws.Cells.Item($row, 8).Value  = 0.29049999999999998
$ws.Cells.Item($row, 9).Value  = 352
$ws.Cells.Item($row, 10).Value = 1.478
$ws.Cells.Item($row, 11).Value = 29
$ws.Cells.Item($row, 12).Value = 1.1020000000000001
$ws.Cells.Item($row, 13).Value = 86
$ws.Cells.Item($row, 14).Value = 1.3759999999999999
$ws.Cells.Item($row, 15).Value = 86
$ws.Cells.Item($row, 16).Value = 0.63439999999999996
$ws.Cells.Item($row, 17).Value = 863
$ws.Cells.Item($row, 18).Value = "Minor bug fixes"

# Match the left-aligned "General" style used by the rest of the data rows.
$ws.Range("A25:R25").HorizontalAlignment = -4131   # xlLeft

# Grow the worksheet table (ListObject) so the new row is part of Table1.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:R25"))

# Move the selection onto the newly added row, like the saved workbook.
$ws.Range("R25").Select()
